$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column E: estado-de-la-informacion -> switch from dimension to measure
$ws.Range("E2").Value = "iaest-measure:estado-de-la-informacion"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("E5").Clear()

# Column J: municipio-nombre -> switch from measure to dimension (refArea)
$ws.Range("J2").Value = "sdmx-dimension:refArea"
$ws.Range("J3").Value = "dim"
$ws.Range("J4").Value = "URI-Municipio"

# Column L: tipo-de-presupuesto -> switch from dimension to measure
$ws.Range("L2").Value = "iaest-measure:tipo-de-presupuesto"
$ws.Range("L3").Value = "medida"
$ws.Range("L4").Value = "xsd:int"
$ws.Range("L5").Clear()
